# =========================================================================
# edit.ps1 - Apply "feat: add 2022-Q4 data" change
#
# Summary of change:
#  1. Insert a new worksheet named "2022-Q4" right after "总计" and before
#     "2022-Q3" (so sheet order becomes 总计, 2022-Q4, 2022-Q3, 2022-Q2, 2022-Q1).
#     The new sheet holds the per-fund holding breakdown for 2022-Q4 (17 rows:
#     1 header row + 16 data rows), formatted the same way as the other
#     quarterly sheets.
#  2. Update the "总计" (summary) sheet: insert a new row 2 with the 2022-Q4
#     totals (持有数量=16, 持有市值=3.04) and keep the existing 2022-Q3 /
#     2022-Q2 / 2022-Q1 rows, shifted down by one row with a corrected index
#     column (A).
# =========================================================================

$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------------
# Step 1: Build the new "2022-Q4" worksheet by cloning the "2022-Q3" sheet.
# This preserves column widths, header styling (bold/centered/bordered,
# style index shared with other quarter sheets) and per-cell formatting,
# so we only need to overwrite the cell values afterwards.
# -------------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("2022-Q3")
$srcSheet.Copy($wb.Worksheets.Item(2))
$ws = $wb.Worksheets.Item(2)
$ws.Name = "2022-Q4"

# The source sheet (2022-Q3) has 24 data rows (rows 2-25); 2022-Q4 only has
# 16 data rows (rows 2-17), so remove the extra trailing rows.
$ws.Rows("18:25").Delete()

# -------------------------------------------------------------------------
# Step 2: Fill in the 2022-Q4 per-fund data.
# Columns B (基金代码), D (基金规模), E (股票总仓位), F (仓位占比) and
# G (持有市值) must stay as TEXT (matching the source data, which stores
# these look-like-numbers values as strings, e.g. fund codes with leading
# zeros such as "010588"). We force the text number format before writing
# the values so the COM layer does not silently coerce them to numbers.
# -------------------------------------------------------------------------

$ws.Range("B2:B17").NumberFormat = "@"
$ws.Range("D2:D17").NumberFormat = "@"
$ws.Range("E2:E17").NumberFormat = "@"
$ws.Range("F2:F17").NumberFormat = "@"
$ws.Range("G2:G17").NumberFormat = "@"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "010588"
$ws.Range("C2").Value = "鹏扬先进制造混合C"
$ws.Range("D2").Value = "8.02"
$ws.Range("E2").Value = "94.05"
$ws.Range("F2").Value = "6.03"
$ws.Range("G2").Value = "0.4836"
$ws.Range("H2").Value = 6

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "010587"
$ws.Range("C3").Value = "鹏扬先进制造混合A"
$ws.Range("D3").Value = "7.90"
$ws.Range("E3").Value = "94.05"
$ws.Range("F3").Value = "6.03"
$ws.Range("G3").Value = "0.4764"
$ws.Range("H3").Value = 6

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "014294"
$ws.Range("C4").Value = "南方北交所精选两年定开混合"
$ws.Range("D4").Value = "4.05"
$ws.Range("E4").Value = "90.08"
$ws.Range("F4").Value = "7.86"
$ws.Range("G4").Value = "0.3183"
$ws.Range("H4").Value = 2

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "005352"
$ws.Range("C5").Value = "鹏扬景泰成长混合A"
$ws.Range("D5").Value = "4.64"
$ws.Range("E5").Value = "94.18"
$ws.Range("F5").Value = "6.66"
$ws.Range("G5").Value = "0.3090"
$ws.Range("H5").Value = 5

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "014203"
$ws.Range("C6").Value = "鹏扬产业趋势一年持有混合A"
$ws.Range("D6").Value = "4.96"
$ws.Range("E6").Value = "94.19"
$ws.Range("F6").Value = "5.70"
$ws.Range("G6").Value = "0.2827"
$ws.Range("H6").Value = 7

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "014271"
$ws.Range("C7").Value = "大成北交所两年定开混合A"
$ws.Range("D7").Value = "3.24"
$ws.Range("E7").Value = "68.93"
$ws.Range("F7").Value = "8.51"
$ws.Range("G7").Value = "0.2757"
$ws.Range("H7").Value = 2

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "014283"
$ws.Range("C8").Value = "华夏北交所创新中小企业精选两年定开混合"
$ws.Range("D8").Value = "3.27"
$ws.Range("E8").Value = "90.95"
$ws.Range("F8").Value = "7.94"
$ws.Range("G8").Value = "0.2596"
$ws.Range("H8").Value = 3

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "014273"
$ws.Range("C9").Value = "广发北交所精选两年定开混合A"
$ws.Range("D9").Value = "3.23"
$ws.Range("E9").Value = "83.79"
$ws.Range("F9").Value = "5.36"
$ws.Range("G9").Value = "0.1731"
$ws.Range("H9").Value = 6

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "014269"
$ws.Range("C10").Value = "嘉实北交所精选两年定期混合A"
$ws.Range("D10").Value = "2.65"
$ws.Range("E10").Value = "94.48"
$ws.Range("F10").Value = "5.89"
$ws.Range("G10").Value = "0.1561"
$ws.Range("H10").Value = 7

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "005353"
$ws.Range("C11").Value = "鹏扬景泰成长混合C"
$ws.Range("D11").Value = "1.26"
$ws.Range("E11").Value = "94.18"
$ws.Range("F11").Value = "6.66"
$ws.Range("G11").Value = "0.0839"
$ws.Range("H11").Value = 5

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "014272"
$ws.Range("C12").Value = "大成北交所两年定开混合C"
$ws.Range("D12").Value = "0.77"
$ws.Range("E12").Value = "68.93"
$ws.Range("F12").Value = "8.51"
$ws.Range("G12").Value = "0.0655"
$ws.Range("H12").Value = 2

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "016307"
$ws.Range("C13").Value = "景顺长城北交所精选两年定开混合A"
$ws.Range("D13").Value = "1.83"
$ws.Range("E13").Value = "43.56"
$ws.Range("F13").Value = "3.00"
$ws.Range("G13").Value = "0.0549"
$ws.Range("H13").Value = 2

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "014274"
$ws.Range("C14").Value = "广发北交所精选两年定开混合C"
$ws.Range("D14").Value = "0.81"
$ws.Range("E14").Value = "83.79"
$ws.Range("F14").Value = "5.36"
$ws.Range("G14").Value = "0.0434"
$ws.Range("H14").Value = 6

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "014270"
$ws.Range("C15").Value = "嘉实北交所精选两年定期混合C"
$ws.Range("D15").Value = "0.52"
$ws.Range("E15").Value = "94.48"
$ws.Range("F15").Value = "5.89"
$ws.Range("G15").Value = "0.0306"
$ws.Range("H15").Value = 7

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "014204"
$ws.Range("C16").Value = "鹏扬产业趋势一年持有混合C"
$ws.Range("D16").Value = "0.25"
$ws.Range("E16").Value = "94.19"
$ws.Range("F16").Value = "5.70"
$ws.Range("G16").Value = "0.0142"
$ws.Range("H16").Value = 7

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "016308"
$ws.Range("C17").Value = "景顺长城北交所精选两年定开混合C"
$ws.Range("D17").Value = "0.27"
$ws.Range("E17").Value = "43.56"
$ws.Range("F17").Value = "3.00"
$ws.Range("G17").Value = "0.0081"
$ws.Range("H17").Value = 2


# -------------------------------------------------------------------------
# Step 3: Update the "总计" (summary) sheet to add the 2022-Q4 row.
# Insert a blank row at row 2 (this shifts 2022-Q3/Q2/Q1 down one row,
# carrying their existing formatting/styles with them), then give the new
# row 2 the same formatting as the data rows (style carried by column A,
# e.g. bold/centered/bordered) before writing the 2022-Q4 totals.
# -------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows("2:2").Insert()
$summary.Range("A2:D2").ClearFormats()

# Re-apply the formatted style (shared with the other index cells in column A)
# to the new A2 cell by copying the format from A3 (the row pushed down,
# which still carries the original style).
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 16
$summary.Range("D2").Value = 3.04

# Fix up the index column (A) for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

# Restore the originally active sheet (the last quarter sheet, 2022-Q1) as
# the selected tab, matching the original workbook's view state.
$wb.Worksheets.Item("2022-Q1").Select()
